$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.041.82"
$ws.Range("D3").Value = "3.189.04"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "535.75"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "145.92"
$ws.Range("E6").Value = "  +4.01%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "0.113"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "3.735.98"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "25.92"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "60.053.83"
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").Value = "3.233.22"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").Value = "6.28"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "8.21"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "368.74"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "0.522"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "69.47"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").Value = "0.169"
$ws.Range("D26").Value = "8.66"
$ws.Range("E26").Value = "  +4.46%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "0.0₃0871"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "22.38"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "6.10"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "5.29"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").Value = "6.57"
$ws.Range("E34").Value = "  +4.54%  "
$ws.Range("D35").Value = "156.27"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "2.812.90"
$ws.Range("E37").Value = "  +6.67%  "
$ws.Range("D38").Value = "26.08"
$ws.Range("E38").Value = "  +3.34%  "
$ws.Range("D39").Value = "0.0704"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0298"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").Value = "39.75"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "3.228.18"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "0.986"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "6.16"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").Value = "20.76"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("E50").Value = "  +4.38%  "
$ws.Range("E51").Value = "  -0.09%  "
